$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Updated "Riders" (column C) and "Average" (column D) values
# New Madigan bike hours

$ws.Range("C2").Value = 284
$ws.Range("D2").Value = 269.62

$ws.Range("C3").Value = 227
$ws.Range("D3").Value = 228.62

$ws.Range("C4").Value = 202
$ws.Range("D4").Value = 197.75

$ws.Range("C5").Value = 201
$ws.Range("D5").Value = 224.86

$ws.Range("C6").Value = 92
$ws.Range("D6").Value = 211.25

$ws.Range("C7").Value = 42
$ws.Range("D7").Value = 122.44
